# EEG and nirs processing 1-23
#
# Row 4 of Sheet1 ("channelsremoved.xlsx") previously listed subject
# NDARVX375BR6 with a removed-channel count of 3. That record is replaced
# with subject NDARAZC45TW3, whose removed-channel count is 0.
#
# (Columns A/B were also nudged slightly wider/narrower by Excel's
# recompute of the stored column width when this row was edited; we
# reproduce that here as closely as the host allows.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "NDARAZC45TW3"
$ws.Range("B4").Value = 0

$ws.Columns.Item(1).ColumnWidth = 14.665
$ws.Columns.Item(2).ColumnWidth = 1.33
